$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.159246203565411
$ws.Range("C2").Value = 4.54990462803945
$ws.Range("D2").Value = 8.979792736928573
$ws.Range("E2").Value = 13.7741381894298
$ws.Range("F2").Value = 34.5532575498221
$ws.Range("I2").Value = 24.06312510769077
$ws.Range("J2").Value = 10.26822494831437
$ws.Range("K2").Value = 9.454984062778024
$ws.Range("M2").Value = 14.97994496646656
$ws.Range("O2").Value = 26.20885599705252

$ws.Range("B3").Value = 8.884693243138791
$ws.Range("C3").Value = 4.367182745457312
$ws.Range("D3").Value = 8.945177260698756
$ws.Range("E3").Value = 13.7707508465181
$ws.Range("F3").Value = 34.63598502883697
$ws.Range("I3").Value = 24.16142677615179
$ws.Range("J3").Value = 10.28965163120326
$ws.Range("K3").Value = 9.278017926341532
$ws.Range("M3").Value = 14.91348651162719
$ws.Range("O3").Value = 26.29806613703064

$ws.Range("B4").Value = 8.712933884357627
$ws.Range("C4").Value = 4.251398255476383
$ws.Range("D4").Value = 8.925247705600649
$ws.Range("E4").Value = 13.77098270969403
$ws.Range("F4").Value = 34.69421619270287
$ws.Range("I4").Value = 24.22599673055972
$ws.Range("J4").Value = 10.30397794892022
$ws.Range("K4").Value = 9.16887682060614
$ws.Range("M4").Value = 14.87470806016475
$ws.Range("O4").Value = 26.35785915870184

$ws.Range("B5").Value = 8.642247163269843
$ws.Range("C5").Value = 4.203379289785845
$ws.Range("D5").Value = 8.917465203559171
$ws.Range("E5").Value = 13.77166008228225
$ws.Range("F5").Value = 34.71981260068038
$ws.Range("I5").Value = 24.25336866349678
$ws.Range("J5").Value = 10.3101106098418
$ws.Range("K5").Value = 9.124336355166767
$ws.Range("M5").Value = 14.85942676693023
$ws.Range("O5").Value = 26.38348551760718

$ws.Range("B6").Value = 8.630471080607144
$ws.Range("C6").Value = 4.195357425596014
$ws.Range("D6").Value = 8.916193562360531
$ws.Range("E6").Value = 13.77180780581888
$ws.Range("F6").Value = 34.72417551263626
$ws.Range("I6").Value = 24.25797772076467
$ws.Range("J6").Value = 10.31114673359718
$ws.Range("K6").Value = 9.116938231947962
$ws.Range("M6").Value = 14.85692112462055
$ws.Range("O6").Value = 26.38781682421463

$ws.Range("B7").Value = 8.711983239181253
$ws.Range("C7").Value = 4.250753944088377
$ws.Range("D7").Value = 8.925141368112262
$ws.Range("E7").Value = 13.77098948297777
$ws.Range("F7").Value = 34.69455384088159
$ws.Range("I7").Value = 24.22636158963592
$ws.Range("J7").Value = 10.30405946302898
$ws.Range("K7").Value = 9.168276317673115
$ws.Range("M7").Value = 14.87449984612213
$ws.Range("O7").Value = 26.35819966420071

$ws.Range("B8").Value = 9.065302719393323
$ws.Range("C8").Value = 4.487685805201154
$ws.Range("D8").Value = 8.967586068079518
$ws.Range("E8").Value = 13.77249146527814
$ws.Range("F8").Value = 34.58023681273117
$ws.Range("I8").Value = 24.09614520729608
$ws.Range("J8").Value = 10.27537016654434
$ws.Range("K8").Value = 9.394100161046397
$ws.Range("M8").Value = 14.95661637013855
$ws.Range("O8").Value = 26.23857334968435

$ws.Range("B9").Value = 9.728600983315742
$ws.Range("C9").Value = 4.921230539540601
$ws.Range("D9").Value = 9.061054958811805
$ws.Range("E9").Value = 13.79370005099422
$ws.Range("F9").Value = 34.41519295800671
$ws.Range("I9").Value = 23.874217346683
$ws.Range("J9").Value = 10.2283825183283
$ws.Range("K9").Value = 9.830679906132623
$ws.Range("M9").Value = 15.13320433210077
$ws.Range("O9").Value = 26.0438653249317

$ws.Range("B10").Value = 10.1926966107678
$ws.Range("C10").Value = 5.217870362417732
$ws.Range("D10").Value = 9.135577970303672
$ws.Range("E10").Value = 13.82029583442825
$ws.Range("F10").Value = 34.33014180081312
$ws.Range("I10").Value = 23.7315484329819
$ws.Range("J10").Value = 10.19949485138645
$ws.Range("K10").Value = 10.14449992273156
$ws.Range("M10").Value = 15.27168703527884
$ws.Range("O10").Value = 25.92521299511394

$ws.Range("B11").Value = 10.39783829902839
$ws.Range("C11").Value = 5.347575969763151
$ws.Range("D11").Value = 9.1706620294134
$ws.Range("E11").Value = 13.8347555226266
$ws.Range("F11").Value = 34.29933812877429
$ws.Range("I11").Value = 23.6710718775938
$ws.Range("J11").Value = 10.18757254545453
$ws.Range("K11").Value = 10.28514281844965
$ws.Range("M11").Value = 15.33641729418297
$ws.Range("O11").Value = 25.87655130245835

$ws.Range("B12").Value = 10.47459216673368
$ws.Range("C12").Value = 5.39590552417813
$ws.Range("D12").Value = 9.184109639700061
$ws.Range("E12").Value = 13.84056745439139
$ws.Range("F12").Value = 34.28880910766716
$ws.Range("I12").Value = 23.64880749698161
$ws.Range("J12").Value = 10.18323282412617
$ws.Range("K12").Value = 10.33804848540862
$ws.Range("M12").Value = 15.36116346214243
$ws.Range("O12").Value = 25.85888997881384

$ws.Range("B13").Value = 10.45810418435092
$ws.Range("C13").Value = 5.385532388001024
$ws.Range("D13").Value = 9.181206380977283
$ws.Range("E13").Value = 13.83930084443041
$ws.Range("F13").Value = 34.29102619127598
$ws.Range("I13").Value = 23.65357419927917
$ws.Range("J13").Value = 10.18415968209156
$ws.Range("K13").Value = 10.32667073540224
$ws.Range("M13").Value = 15.35582375885723
$ws.Range("O13").Value = 25.86265957547581

$ws.Range("B14").Value = 10.40417187220264
$ws.Range("C14").Value = 5.351568024232155
$ws.Range("D14").Value = 9.171765172142909
$ws.Range("E14").Value = 13.8352269516636
$ws.Range("F14").Value = 34.29844913198598
$ws.Range("I14").Value = 23.66922740819256
$ws.Range("J14").Value = 10.18721200833682
$ws.Range("K14").Value = 10.28950273729182
$ws.Range("M14").Value = 15.33844856977912
$ws.Range("O14").Value = 25.87508293749696

$ws.Range("B15").Value = 10.37101387993498
$ws.Range("C15").Value = 5.330660435883114
$ws.Range("D15").Value = 9.166003016889601
$ws.Range("E15").Value = 13.83277528031809
$ws.Range("F15").Value = 34.30314383379179
$ws.Range("I15").Value = 23.67889839565694
$ws.Range("J15").Value = 10.18910442795898
$ws.Range("K15").Value = 10.2666889337177
$ws.Range("M15").Value = 15.32783581384108
$ws.Range("O15").Value = 25.88279238442432

$ws.Range("B16").Value = 10.17916374613208
$ws.Range("C16").Value = 5.209285384771944
$ws.Range("D16").Value = 9.133308266035071
$ws.Range("E16").Value = 13.81939810738185
$ws.Range("F16").Value = 34.33231373594193
$ws.Range("I16").Value = 23.73558979385117
$ws.Range("J16").Value = 10.20029850530877
$ws.Range("K16").Value = 10.1352617980273
$ws.Range("M16").Value = 15.26749036668542
$ws.Range("O16").Value = 25.92850022421341

$ws.Range("B17").Value = 10.05988726923603
$ws.Range("C17").Value = 5.133458499967137
$ws.Range("D17").Value = 9.113548447938864
$ws.Range("E17").Value = 13.81179422163287
$ws.Range("F17").Value = 34.35222972166972
$ws.Range("I17").Value = 23.77150168657463
$ws.Range("J17").Value = 10.20747769273155
$ws.Range("K17").Value = 10.05405828502698
$ws.Range("M17").Value = 15.23090382850455
$ws.Range("O17").Value = 25.95790272036283

$ws.Range("B18").Value = 9.99072360511833
$ws.Range("C18").Value = 5.089353613322128
$ws.Range("D18").Value = 9.102294989456759
$ws.Range("E18").Value = 13.80764311553945
$ws.Range("F18").Value = 34.36442711909768
$ws.Range("I18").Value = 23.79257365205581
$ws.Range("J18").Value = 10.21172171088729
$ws.Range("K18").Value = 10.00715579286723
$ws.Range("M18").Value = 15.21002433841843
$ws.Range("O18").Value = 25.97531432617534

$ws.Range("B19").Value = 9.967212199359867
$ws.Range("C19").Value = 5.074337199918032
$ws.Range("D19").Value = 9.098504217707266
$ws.Range("E19").Value = 13.80627591998603
$ws.Range("F19").Value = 34.36868437485314
$ws.Range("I19").Value = 23.79977975666861
$ws.Range("J19").Value = 10.21317837725175
$ws.Range("K19").Value = 9.991243211065289
$ws.Range("M19").Value = 15.20298354346154
$ws.Range("O19").Value = 25.98129543174264

$ws.Range("B20").Value = 10.07264283174371
$ws.Range("C20").Value = 5.141581507927394
$ws.Range("D20").Value = 9.115640394483403
$ws.Range("E20").Value = 13.81258066780867
$ws.Range("F20").Value = 34.35003279808927
$ws.Range("I20").Value = 23.76763570492171
$ws.Range("J20").Value = 10.20670158281121
$ws.Range("K20").Value = 10.06272322258928
$ws.Range("M20").Value = 15.2347816504653
$ws.Range("O20").Value = 25.95472100696027

$ws.Range("B21").Value = 10.42003882409992
$ws.Range("C21").Value = 5.36156579080801
$ws.Range("D21").Value = 9.174533952076425
$ws.Range("E21").Value = 13.83641445054428
$ws.Range("F21").Value = 34.29623800235532
$ws.Range("I21").Value = 23.6646123945509
$ws.Range("J21").Value = 10.18631071894319
$ws.Range("K21").Value = 10.30042982588901
$ws.Range("M21").Value = 15.34354584919392
$ws.Range("O21").Value = 25.87141309816776

$ws.Range("B22").Value = 10.64164109473526
$ws.Range("C22").Value = 5.500737733051988
$ws.Range("D22").Value = 9.213964863446456
$ws.Range("E22").Value = 13.85395011882942
$ws.Range("F22").Value = 34.26769992144026
$ws.Range("I22").Value = 23.60099262857425
$ws.Range("J22").Value = 10.17400399796279
$ws.Range("K22").Value = 10.45370821652949
$ws.Range("M22").Value = 15.41598754595225
$ws.Range("O22").Value = 25.82143098646861

$ws.Range("B23").Value = 10.52388668457104
$ws.Range("C23").Value = 5.426889988511756
$ws.Range("D23").Value = 9.192836517418886
$ws.Range("E23").Value = 13.84441287797079
$ws.Range("F23").Value = 34.28232509694371
$ws.Range("I23").Value = 23.63460785467286
$ws.Range("J23").Value = 10.18047909784671
$ws.Range("K23").Value = 10.37210598472508
$ws.Range("M23").Value = 15.37720486285874
$ws.Range("O23").Value = 25.84769834652509

$ws.Range("B24").Value = 10.06687787440904
$ws.Range("C24").Value = 5.137910685052423
$ws.Range("D24").Value = 9.114694292415193
$ws.Range("E24").Value = 13.81222442852621
$ws.Range("F24").Value = 34.35102369915734
$ws.Range("I24").Value = 23.76938218910158
$ws.Range("J24").Value = 10.20705209886344
$ws.Range("K24").Value = 10.05880647809282
$ws.Range("M24").Value = 15.23302800445891
$ws.Range("O24").Value = 25.95615787818877

$ws.Range("B25").Value = 9.552894572744057
$ws.Range("C25").Value = 4.807601966559889
$ws.Range("D25").Value = 9.034712222688835
$ws.Range("E25").Value = 13.78601866735588
$ws.Range("F25").Value = 34.45349340751782
$ws.Range("I25").Value = 23.93067593393564
$ws.Range("J25").Value = 10.24010316240562
$ws.Range("K25").Value = 9.713574018842991
$ws.Range("M25").Value = 15.08384152379499
$ws.Range("O25").Value = 26.0438653249317
